# Update cryptocurrency price/volume data to latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.514.24"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.44%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.869.57"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -0.53%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.95%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "315.19"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("E6").Value = "  -0.57%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5054"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -1.57%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3892"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.81%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08344"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.42%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "41.74"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -0.36%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.100"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -1.83%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.193"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.37%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.870.69"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -0.22%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.35"
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.225"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -0.35%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -0.89%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001099"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -0.66%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "90.86"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.38%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06692"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -0.59%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.63"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("E21").Value = "  -0.66%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.904"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.66%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.549.11"
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.03"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.17%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.232"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -1.09%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.085.78"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.24%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "161.69"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.54%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "20.61"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -1.04%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.334"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -4.57%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "125.54"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -0.56%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.1040"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("E32").Value = "  -1.04%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.766"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -2.25%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.603"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.82%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.02445"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.15%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.06518"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -0.21%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.2155"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -1.63%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "8.830"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -4.45%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.041"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.21%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.252"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.184"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.49%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6404"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -1.13%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "11.06"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.15%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -0.53%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.6009"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.99%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "12.98"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.51%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.691"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("E49").Value = "  -0.87%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "121.63"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.34%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.177"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -8.24%  "
